# Insert a new data row at row 46 (pushing existing rows 46-101 down to 47-102)
# and populate it with a new "Alcachofa" price record for the Macroferia Regional
# de Talca market (Maule), matching the weekly fruit/vegetable update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A46").EntireRow.Insert()

$ws.Range("A46").Value = 5
$ws.Range("B46").Value = "Macroferia Regional de Talca"
$ws.Range("C46").Value = "Maule"
$ws.Range("D46").Value = (Get-Date -Year 2022 -Month 8 -Day 30 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E46").Value = 7
$ws.Range("F46").Value = 100112013
$ws.Range("G46").Value = "Alcachofa"
$ws.Range("H46").Value = "Madrigal"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 200
$ws.Range("K46").Value = 13000
$ws.Range("L46").Value = 13000
$ws.Range("M46").Value = 13000
$ws.Range("N46").Value = "`$/caja 40 unidades"
$ws.Range("O46").Value = "Provincia del Elquí"
$ws.Range("P46").Value = 325
$ws.Range("Q46").Value = 40
$ws.Range("R46").Value = "Hortaliza"
